$d = $word.ActiveDocument

# 1) "Aplikacja ktorą" -> "Aplikacja, ktorą" (insert comma after the first word)
$null = $d.Content.Find.Execute(
    "Aplikacja którą", $true, $false, $false, $false, $false,
    $true, 1, $false, "Aplikacja, którą", 2)

# 2) Rework the description of what the application manages.
$null = $d.Content.Find.Execute(
    "zarządzać sprzedażami firmy sprzedającej ekskluzywne kosmetyki: podkłady i korektory. Aplikacja służyć będzie nieznającym",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "zarządzać stanem magazynowym oraz sprzedażą sklepu z kosmetykami (podkładami i korektorami) luksusowych marek. Będzie służyć nieznającym",
    2)

# 3) "Aplikacja dla podstawowego" -> "Dla podstawowego"
$null = $d.Content.Find.Execute(
    " Aplikacja dla podstawowego", $true, $false, $false, $false, $false,
    $true, 1, $false, " Dla podstawowego", 2)
